$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") rows 2-38, replacing the old Strike# values
$newValues = @{
    2  = 5
    3  = 10
    4  = 6
    5  = 8
    6  = 7
    7  = 7
    8  = 2
    9  = 4
    10 = 8
    11 = 3
    12 = 6
    13 = 6
    14 = 3
    15 = 9
    16 = 4
    17 = 6
    18 = 8
    19 = 8
    20 = 4
    21 = 5
    22 = 2
    23 = 6
    24 = 7
    25 = 5
    26 = 8
    27 = 3
    28 = 7
    29 = 4
    30 = 9
    31 = 3
    32 = 6
    33 = 5
    34 = 7
    35 = 4
    36 = 5
    37 = 2
    38 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
